$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. De-select the "naive" tab (it will no longer be the active sheet once we
#    add + activate the new sheets below; Excel auto-clears tabSelected on
#    every other sheet when a new one becomes active).
# ---------------------------------------------------------------------------
$wsSibreg = $wb.Worksheets.Item("sibregsimple")
$wsNaive  = $wb.Worksheets.Item("naive")

# ---------------------------------------------------------------------------
# 2. Add "nass_sibreg" sheet right after "naive" -- same layout/columns as
#    "sibregsimple" (runyear / predicted return / p25 / p75) but only the
#    runyear column is populated; the rest are blank awaiting new data.
# ---------------------------------------------------------------------------
$wsNassSibreg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNaive)
$wsNassSibreg.Name = "nass_sibreg"

$wsNassSibreg.Range("A1").Value = "runyear"
$wsNassSibreg.Range("B1").Value = "predicted return"
$wsNassSibreg.Range("C1").Value = "p25"
$wsNassSibreg.Range("D1").Value = "p75"

for ($i = 0; $i -lt 14; $i++) {
    $row = $i + 2
    $wsNassSibreg.Cells.Item($row, 1).Value = 2012 + $i
    $wsNassSibreg.Cells.Item($row, 1).NumberFormat = "0"
    $wsNassSibreg.Cells.Item($row, 2).NumberFormat = "0.00"
    $wsNassSibreg.Cells.Item($row, 3).NumberFormat = "0"
    $wsNassSibreg.Cells.Item($row, 4).NumberFormat = "0"
}

$loSibreg = $wsNassSibreg.ListObjects.Add(1, $wsNassSibreg.Range("A1:D15"), [System.Reflection.Missing]::Value, 1)
$loSibreg.Name = "Table14"
$loSibreg.TableStyle = "TableStyleLight1"

$wsNassSibreg.Range("A1:D15").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. Add "nass_naive" sheet right after "nass_sibreg" -- same layout as
#    "naive", again only runyear populated.
# ---------------------------------------------------------------------------
$wsNassNaive = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNassSibreg)
$wsNassNaive.Name = "nass_naive"

$wsNassNaive.Range("A1").Value = "runyear"
$wsNassNaive.Range("B1").Value = "predicted return"
$wsNassNaive.Range("C1").Value = "p25"
$wsNassNaive.Range("D1").Value = "p75"

for ($i = 0; $i -lt 14; $i++) {
    $row = $i + 2
    $wsNassNaive.Cells.Item($row, 1).Value = 2012 + $i
    $wsNassNaive.Cells.Item($row, 1).NumberFormat = "0"
    $wsNassNaive.Cells.Item($row, 2).NumberFormat = "0.00"
    $wsNassNaive.Cells.Item($row, 3).NumberFormat = "0"
    $wsNassNaive.Cells.Item($row, 4).NumberFormat = "0"
}

$loNaive = $wsNassNaive.ListObjects.Add(1, $wsNassNaive.Range("A1:D15"), [System.Reflection.Missing]::Value, 1)
$loNaive.Name = "Table145"
$loNaive.TableStyle = "TableStyleLight1"

$wsNassNaive.Range("C23").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Make "nass_sibreg" the active tab (matches activeTab="2" / tabSelected
#    moving to the 3rd sheet).
# ---------------------------------------------------------------------------
$wsNassSibreg.Activate() | Out-Null
$wsNassSibreg.Range("A1:D15").Select() | Out-Null
